# Update cryptocurrency Price (D) and Volume(1h) (E) columns with refreshed
# GitHub Actions scrape values, row by row (rows 2-51 of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.147.31"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.839.88"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.92"
$ws.Range("E5").Value = "  -2.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6860"
$ws.Range("E6").Value = "  -2.62%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2992"
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07424"
$ws.Range("E9").Value = "  -4.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.18"
$ws.Range("E10").Value = "  -2.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07644"
$ws.Range("E11").Value = "  -2.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.843.18"
$ws.Range("E12").Value = "  -0.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.049"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6807"
$ws.Range("E14").Value = "  -1.90%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "87.43"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.145"
$ws.Range("E16").Value = "  -7.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.153.73"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008164"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.085.43"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "229.52"
$ws.Range("E20").Value = "  -5.87%  "
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.365"
$ws.Range("E23").Value = "  -2.82%  "
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.34"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1441"
$ws.Range("E26").Value = "  -5.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.749"
$ws.Range("E27").Value = "  -1.84%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.08"
$ws.Range("E28").Value = "  -1.75%  "
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.263"
$ws.Range("E30").Value = "  +0.24%  "
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05271"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7567"
$ws.Range("E34").Value = "  -4.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.854"
$ws.Range("E35").Value = "  -4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.133"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.294.87"
$ws.Range("E38").Value = "  -3.12%  "
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.719"
$ws.Range("E40").Value = "  -0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9381"
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.954"
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "104.79"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  +0.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.988.03"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "64.84"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.489"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.766"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07401"
$ws.Range("E51").Value = "  +17.44%  "
